$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "URWpD674"
$ws.Range("B2").Value = 231004199
$ws.Range("C2").Value = "dquwzpx93"
$ws.Range("D2").Value = "AM!j`$6w5"
$ws.Range("F2").Value = "hnhwTiVw"
$ws.Range("G2").Value = "cDtg"

# Row 3
$ws.Range("A3").Value = "JUNzx491"
$ws.Range("B3").Value = 231004198
$ws.Range("C3").Value = "ojlgjzu53"
$ws.Range("D3").Value = "b6dY%4U`$"
$ws.Range("F3").Value = "RLBqYUNo"
$ws.Range("G3").Value = "Gczo"
